{"js": "// Apply the edits described by the commit \"upload do an 1 + sua 1 so printf\":\n//  1) Table row \"1c\": \"3-change: ...\" -> \"3-change info: ...\"\n//  2) Table row \"1d\": \"... ? : X\u00f3a th\u00f4ng tin \u0111\u1ed9c gi\u1ea3 qua t\u00ean\"\n//        -> \"... ? : T\u00ecm th\u00f4ng tin c\u1ee7a \u0111\u1ed9c gi\u1ea3 b\u1eb1ng t\u00ean + x\u00e1c nh\u1eadn xo\u00e1 \"\n//  3) Table row \"1f\": \"6-search book borrowed by name and delete ?: xem th\u00f4ng tin\n//        s\u00e1ch \u0111\u00e3 m\u01b0\u1ee3n c\u1ee7a \u0111\u1ed9c gi\u1ea3 (s\u1eed d\u1ee5ng chung v\u1edbi 1d, ch\u1ecdn kh\u00f4ng x\u00f3a)\"\n//        -> \"2-search and get info by name : xem th\u00f4ng tin c\u1ee7a \u0111\u1ed9c gi\u1ea3 \"\n//  4) Bullet list item \"T\u00ecm ki\u1ebfm s\u00e1ch m\u01b0\u1ee3n theo t\u00ean v\u00e0 x\u00f3a (search book borrowed\n//        by name and delete)\" -> \"T\u00ecm ki\u1ebfm th\u00f4ng tin theo t\u00ean v\u00e0 x\u00f3a (search book\n//        borrowed by name and delete)\", and a new bullet item\n//        \"T\u00ecm ki\u1ebfm theo t\u00ean (search and get info by name)\" is added right after it.\n\nconst body = context.document.body;\n\n// --- 1) \"3-change: ch\u1ec9nh s\u1eeda th\u00f4ng tin 1 \u0111\u1ed9c gi\u1ea3\" -----------------------------\nconst r1 = body.search(\"3-change: ch\u1ec9nh s\u1eeda th\u00f4ng tin 1 \u0111\u1ed9c gi\u1ea3\", { matchCase: true });\nr1.load(\"items\");\nawait context.sync();\nif (r1.items.length > 0) {\n  r1.items[0].insertText(\"3-change info: ch\u1ec9nh s\u1eeda th\u00f4ng tin 1 \u0111\u1ed9c gi\u1ea3\", \"Replace\");\n}\n\n// --- 2) \"...delete ? : X\u00f3a th\u00f4ng tin \u0111\u1ed9c gi\u1ea3 qua t\u00ean\" -------------------------\nconst r2 = body.search(\n  \"6-search book borrowed by name and delete ? : X\u00f3a th\u00f4ng tin \u0111\u1ed9c gi\u1ea3 qua t\u00ean\",\n  { matchCase: true }\n);\nr2.load(\"items\");\nawait context.sync();\nif (r2.items.length > 0) {\n  r2.items[0].insertText(\n    \"6-search book borrowed by name and delete ? : T\u00ecm th\u00f4ng tin c\u1ee7a \u0111\u1ed9c gi\u1ea3 b\u1eb1ng t\u00ean + x\u00e1c nh\u1eadn xo\u00e1 \",\n    \"Replace\"\n  );\n}\n\n// --- 3) \"...delete ?: xem th\u00f4ng tin s\u00e1ch \u0111\u00e3 m\u01b0\u1ee3n ...\" -------------------------\nconst r3 = body.search(\n  \"6-search book borrowed by name and delete ?: xem th\u00f4ng tin s\u00e1ch \u0111\u00e3 m\u01b0\u1ee3n c\u1ee7a \u0111\u1ed9c gi\u1ea3 (s\u1eed d\u1ee5ng chung v\u1edbi 1d, ch\u1ecdn kh\u00f4ng x\u00f3a)\",\n  { matchCase: true }\n);\nr3.load(\"items\");\nawait context.sync();\nif (r3.items.length > 0) {\n  r3.items[0].insertText(\n    \"2-search and get info by name : xem th\u00f4ng tin c\u1ee7a \u0111\u1ed9c gi\u1ea3 \",\n    \"Replace\"\n  );\n}\n\n// --- 4) Heading bullet: \"s\u00e1ch m\u01b0\u1ee3n\" -> \"th\u00f4ng tin\", plus a new bullet item ----\nconst r4 = body.search(\n  \"T\u00ecm ki\u1ebfm s\u00e1ch m\u01b0\u1ee3n theo t\u00ean v\u00e0 x\u00f3a (search book borrowed by name and delete)\",\n  { matchCase: true }\n);\nr4.load(\"items\");\nawait context.sync();\nif (r4.items.length > 0) {\n  const headingRange = r4.items[0];\n  headingRange.insertText(\n    \"T\u00ecm ki\u1ebfm th\u00f4ng tin theo t\u00ean v\u00e0 x\u00f3a (search book borrowed by name and delete)\",\n    \"Replace\"\n  );\n  await context.sync();\n\n  const headingPara = headingRange.paragraphs.getFirst();\n  const newPara = headingPara.insertParagraph(\n    \"T\u00ecm ki\u1ebfm theo t\u00ean (search and get info by name)\",\n    \"After\"\n  );\n  // Match the bold run formatting used by the surrounding list items.\n  newPara.font.bold = true;\n}\n\nawait context.sync();\n", "ps1": "# Apply the edits described by the commit \"upload do an 1 + sua 1 so printf\":\n#  1) Table row \"1c\": \"3-change: ...\" -> \"3-change info: ...\"\n#  2) Table row \"1d\": \"... ? : X\u00f3a th\u00f4ng tin \u0111\u1ed9c gi\u1ea3 qua t\u00ean\"\n#        -> \"... ? : T\u00ecm th\u00f4ng tin c\u1ee7a \u0111\u1ed9c gi\u1ea3 b\u1eb1ng t\u00ean + x\u00e1c nh\u1eadn xo\u00e1 \"\n#  3) Table row \"1f\": \"6-search book borrowed by name and delete ?: xem th\u00f4ng tin\n#        s\u00e1ch \u0111\u00e3 m\u01b0\u1ee3n c\u1ee7a \u0111\u1ed9c gi\u1ea3 (s\u1eed d\u1ee5ng chung v\u1edbi 1d, ch\u1ecdn kh\u00f4ng x\u00f3a)\"\n#        -> \"2-search and get info by name : xem th\u00f4ng tin c\u1ee7a \u0111\u1ed9c gi\u1ea3 \"\n#  4) Bullet list item \"T\u00ecm ki\u1ebfm s\u00e1ch m\u01b0\u1ee3n theo t\u00ean v\u00e0 x\u00f3a (search book borrowed\n#        by name and delete)\" -> \"T\u00ecm ki\u1ebfm th\u00f4ng tin theo t\u00ean v\u00e0 x\u00f3a (search book\n#        borrowed by name and delete)\", and a new bullet item\n#        \"T\u00ecm ki\u1ebfm theo t\u00ean (search and get info by name)\" is added right after it.\n\n$d = $word.ActiveDocument\n\n# wdReplaceOne = 1, wdFindWrapFind (wdFindWrap) = 1, wdFindContinue = 1\n$wdReplaceOne = 1\n\n# --- 1) \"3-change: ch\u1ec9nh s\u1eeda th\u00f4ng tin 1 \u0111\u1ed9c gi\u1ea3\" -----------------------------\n$rng1 = $d.Content\n$rng1.Find.Execute(\n    \"3-change: ch\u1ec9nh s\u1eeda th\u00f4ng tin 1 \u0111\u1ed9c gi\u1ea3\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"3-change info: ch\u1ec9nh s\u1eeda th\u00f4ng tin 1 \u0111\u1ed9c gi\u1ea3\",\n    $wdReplaceOne\n) | Out-Null\n\n# --- 2) \"...delete ? : X\u00f3a th\u00f4ng tin \u0111\u1ed9c gi\u1ea3 qua t\u00ean\" -------------------------\n$rng2 = $d.Content\n$rng2.Find.Execute(\n    \"6-search book borrowed by name and delete ? : X\u00f3a th\u00f4ng tin \u0111\u1ed9c gi\u1ea3 qua t\u00ean\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"6-search book borrowed by name and delete ? : T\u00ecm th\u00f4ng tin c\u1ee7a \u0111\u1ed9c gi\u1ea3 b\u1eb1ng t\u00ean + x\u00e1c nh\u1eadn xo\u00e1 \",\n    $wdReplaceOne\n) | Out-Null\n\n# --- 3) \"...delete ?: xem th\u00f4ng tin s\u00e1ch \u0111\u00e3 m\u01b0\u1ee3n ...\" -------------------------\n$rng3 = $d.Content\n$rng3.Find.Execute(\n    \"6-search book borrowed by name and delete ?: xem th\u00f4ng tin s\u00e1ch \u0111\u00e3 m\u01b0\u1ee3n c\u1ee7a \u0111\u1ed9c gi\u1ea3 (s\u1eed d\u1ee5ng chung v\u1edbi 1d, ch\u1ecdn kh\u00f4ng x\u00f3a)\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"2-search and get info by name : xem th\u00f4ng tin c\u1ee7a \u0111\u1ed9c gi\u1ea3 \",\n    $wdReplaceOne\n) | Out-Null\n\n# --- 4) Heading bullet: \"s\u00e1ch m\u01b0\u1ee3n\" -> \"th\u00f4ng tin\", plus a new bullet item ----\n$rng4 = $d.Content\n$rng4.Find.Execute(\n    \"T\u00ecm ki\u1ebfm s\u00e1ch m\u01b0\u1ee3n theo t\u00ean v\u00e0 x\u00f3a (search book borrowed by name and delete)\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"T\u00ecm ki\u1ebfm th\u00f4ng tin theo t\u00ean v\u00e0 x\u00f3a (search book borrowed by name and delete)\",\n    $wdReplaceOne\n) | Out-Null\n\n# Locate the (now renamed) heading paragraph again and append a sibling bullet\n# paragraph right after it with the new text.\n$rng5 = $d.Content\n$rng5.Find.Execute(\"T\u00ecm ki\u1ebfm th\u00f4ng tin theo t\u00ean v\u00e0 x\u00f3a (search book borrowed by name and delete)\") | Out-Null\n$headingPara = $rng5.Paragraphs.First\n$headingRange = $headingPara.Range\n$headingRange.Collapse(0)\n$headingRange.InsertParagraphAfter()\n\n$newItemRange = $headingPara.Range\n$newItemRange.Collapse(0)\n$newItemRange.MoveEnd(1, 1) | Out-Null\n$newItemRange.InsertAfter(\"T\u00ecm ki\u1ebfm theo t\u00ean (search and get info by name)\")\n\nWrite-Output \"done\"\n"}
